$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Clientes" (sheet1) — add Numero/Complemento/Bairro/CEP/Telefone
# columns, split old "Cidade"/"Estado" pair further right, and split
# the old "Rua X, NN" address string into street + number.
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Clientes")

# Header row
$ws1.Range("A1").Value = "Identificador"
$ws1.Range("B1").Value = "Razão Social"
$ws1.Range("C1").Value = "CNPJ"
$ws1.Range("D1").Value = "Endereço"
$ws1.Range("E1").Value = "Número"
$ws1.Range("F1").Value = "Complemento"
$ws1.Range("G1").Value = "Bairro"
$ws1.Range("H1").Value = "Cidade"
$ws1.Range("I1").Value = "Estado"
$ws1.Range("J1").Value = "CEP"
$ws1.Range("K1").Value = "Telefone"

# Row 2 — Eventos & Cia
$ws1.Range("A2").Value = 1
$ws1.Range("B2").Value = "Eventos & Cia"
$ws1.Range("C2").Value = "62.100.755/0001-15"
$ws1.Range("D2").Value = "Rua Presidente Kennedy"
$ws1.Range("E2").Value = 57
$ws1.Range("F2").Value = "Casa 05"
$ws1.Range("G2").Value = "Cônego"
$ws1.Range("H2").Value = "Nova Friburgo"
$ws1.Range("I2").Value = "RJ"
$ws1.Range("J2").Value = 28621000
$ws1.Range("K2").Value = "(22) 2522-5120"

# Row 3 — Petrobras
$ws1.Range("A3").Value = 2
$ws1.Range("B3").Value = "Petrobras - Petróleo Brasileiro SA"
$ws1.Range("C3").Value = "78.865.726/0001-84"
$ws1.Range("D3").Value = "Avenida República do Chile"
$ws1.Range("E3").Value = 65
$ws1.Range("F3").ClearContents()
$ws1.Range("G3").Value = "Centro"
$ws1.Range("H3").Value = "Rio de Janeiro"
$ws1.Range("I3").Value = "RJ"
$ws1.Range("J3").Value = 20031912
$ws1.Range("K3").Value = "(22) 2513-0056"

# ---------------------------------------------------------------
# Sheet "Dados Manjerico" (sheet2) — same column additions.
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Dados Manjerico")

# Header row
$ws2.Range("A1").Value = "Razão Social"
$ws2.Range("B1").Value = "CNPJ"
$ws2.Range("C1").Value = "Endereço"
$ws2.Range("D1").Value = "Número"
$ws2.Range("E1").Value = "Complemento"
$ws2.Range("F1").Value = "Bairro"
$ws2.Range("G1").Value = "Cidade"
$ws2.Range("H1").Value = "Estado"
$ws2.Range("I1").Value = "CEP"
$ws2.Range("J1").Value = "Telefone"

# Row 2 — Manjerico Corp SA
$ws2.Range("A2").Value = "Manjerico Corp SA"
$ws2.Range("B2").Value = "06.828.467/0001-00"
$ws2.Range("C2").Value = "Rua Itacuruçá, 26"
$ws2.Range("D2").Value = 26
$ws2.Range("E2").Value = "Apto 212"
$ws2.Range("F2").Value = "Tijuca"
$ws2.Range("G2").Value = "Rio de Janeiro"
$ws2.Range("H2").Value = "RJ"
$ws2.Range("I2").Value = 20510150
$ws2.Range("J2").Value = "(21) 2135-1448"

$ws2.Range("F2").Select() | Out-Null
$ws1.Range("G4").Select() | Out-Null
